# Update average_county_temperature (I), worst_ashp_cop (N) and best_ashp_cop (O)
# with refreshed NOAA-sourced temperature data and recalculated COP values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row=8;  I=13.76976495726495; N=1.867772014163364; O=2.045100507661769 },
    @{ Row=14; I=21.19907407407406; N=2.015705049109126; O=2.22640738080769 },
    @{ Row=15; I=2.356481481481501; N=1.678525338046114; O=1.817698795724144 },
    @{ Row=20; I=20.68981481481483; N=2.004820578273036; O=2.212959200483225 },
    @{ Row=21; I=21.28240740740739; N=2.017497406510892; O=2.228623569098047 },
    @{ Row=22; I=15.74228395061728; N=1.904889690449167; O=2.090295475371289 },
    @{ Row=23; I=-1.819444444444444; N=1.618523362263702; O=1.746638928617865 },
    @{ Row=25; I=13.46442495126706; N=1.862155209238257; O=2.038278558917324 },
    @{ Row=29; I=19.79629629629628; N=1.98600466835246;  O=2.18975222777657 },
    @{ Row=30; I=19.79629629629628 },
    @{ Row=32; I=21.79166666666666; N=2.028520339740724; O=2.242263395092639 },
    @{ Row=33; I=12.51681286549706; N=1.844936767548521; O=2.017393709936214 },
    @{ Row=34; I=15.74228395061728; N=1.904889690449167; O=2.090295475371289 },
    @{ Row=35; I=15.74228395061728 },
    @{ Row=39; I=13.75752314814816; N=1.867546171126113; O=2.044826120875009 },
    @{ Row=40; I=14.96875 },
    @{ Row=41; I=14.96875;          N=1.890159325210871; O=2.072335994446373 },
    @{ Row=52; I=19.48611111111111; N=1.979555038534245; O=2.181809322722105 },
    @{ Row=53; I=14.47727272727272; N=1.880917929007461; O=2.06108460959076 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("I$r").Value = $u.I
    if ($u.ContainsKey("N")) {
        $ws.Range("N$r").Value = $u.N
    }
    if ($u.ContainsKey("O")) {
        $ws.Range("O$r").Value = $u.O
    }
}
